# Scenario 13 / TC001 - "Upload Incoterms" sheet update.
# The test-case code value in the "TB Incoterms Code" row was bumped
# from S4INTCODET -> S6INTCODET, and the last-saved cell selection
# moved from E25 to C27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MGF060")

# C21 holds the "TB Incoterms Code" test value (D21 is the label).
$ws.Range("C21").Value = "S6INTCODET"

# Leave the cursor where the author last left it before saving.
$ws.Range("C27").Select()
